$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.021.29"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "2.737.03"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.24%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "2.735.61"
$ws.Range("E9").Value = "  +3.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.20%  "
$ws.Range("E11").Value = "  +5.58%  "
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "3.235.06"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").Value = "68.948.04"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "2.729.65"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "376.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("E21").Value = "  +5.12%  "
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("E23").Value = "  +6.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.50%  "
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000106"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "589.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  +5.47%  "
$ws.Range("E33").Value = "  +5.49%  "
$ws.Range("E34").Value = "  +5.83%  "
$ws.Range("E35").Value = "  +4.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "162.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("E41").Value = "  +4.18%  "
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.33%  "
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("D47").Value = "0.0₆0310"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  +5.55%  "
$ws.Range("E50").Value = "  +7.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.607"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.97%  "
Write-Host "Updated cryptos list"
